$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.41111872142219852,
    -0.060725433925082228,
    -0.0089999997117047315,
    -0.011999999922139182,
    -0.0059999997146276129,
    -0.0059999997046276121,
    -0.019999999649272127,
    -0.019999999644886302,
    -0.0059999996934472222,
    -0.0059999996887185603,
    -0.0044999996939800724,
    -0.0059999996870008232,
    -0.0059999996819239954,
    -0.011999999657259508,
    -0.0059999996791475496,
    -0.0059999996780386589,
    -0.0059999996765736086,
    -0.0089999996646294989,
    -0.0089999997205287841,
    -0.031265225203304325,
    -0.0089999997155514322,
    -0.0089999997152503397,
    -0.0089999997092053974,
    -0.041999999577526204,
    -0.04199999957517786,
    -0.0059999997041266795,
    -0.0059999997034032582,
    -0.005999999700828873,
    -0.011999999676191919,
    -0.019999999644922717,
    0.031590929329297524,
    -0.020999999640889833,
    -0.0059999996984299031
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
